# Week2 Day7 session management analysis - mark Day4/Day5/Day6 "완료" (G5:G7)
# cells as checked, using a distinct "checkbox" glyph font (Segoe UI Symbol),
# and move the active selection to F16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Checked-box glyph (☑ U+2611)
$checkGlyph = [char]0x2611

# Set G5 first (this is where the new font gets introduced), then copy its
# formatting onto G6:G7 so every one of them ends up sharing the exact same
# style record instead of each creating its own.
$ws.Range("G5").Value = $checkGlyph
$ws.Range("G5").Font.Name = "Segoe UI Symbol"

$ws.Range("G5").Copy()
$ws.Range("G6:G7").PasteSpecial(-4122)
$ws.Range("G6:G7").Value = $checkGlyph
$excel.CutCopyMode = 0

# Move the selection, matching the saved sheet view state.
$null = $ws.Range("F16").Select()
